# Y4_B2526_General_Surgery_checklist — attendance-app re-upload edit
#
# The log sheet's "Log Date" column (C) is bumped forward by 30 days
# (8/9/2025 -> 9/8/2025, serial 45878 -> 45908) for every logged row,
# and the "Log Time" column (D) is re-formatted to a 12-hour clock with
# an explicit AM/PM marker. Column C also gets an explicit width, the
# rows drop their stale explicit height (falling back to the sheet's
# default row height) and the active selection is reset to the header
# row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

$firstDataRow = 2
$lastDataRow  = 74

$dataRange = $ws.Range("C" + $firstDataRow + ":C" + $lastDataRow)

# 1. Shift every "Log Date" entry forward by 30 days (45878 -> 45908).
$dataRange.Value = "9/8/2025"

# 2. Number formats: date for column C, 12-hour AM/PM time for column D.
#    Re-applied across the header too, matching how Excel stamps a
#    whole column's cells (including the header) once a column format
#    is set.
$ws.Range("C1:C" + $lastDataRow).NumberFormat = "m/d/yyyy"
$ws.Range("D1:D" + $lastDataRow).NumberFormat = "[$-F400]h:mm:ss AM/PM"

# 3. Give column C an explicit width (column D already has one).
$ws.Columns.Item(3).ColumnWidth = 14.63

# 4. Rows no longer carry an explicit 15.75pt height override; let them
#    fall back to the sheet's default row height.
$ws.Range("A1:F" + $lastDataRow).EntireRow.AutoFit()

# 5. Reset the active selection to the header row.
$ws.Rows.Item(1).Select()
